$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$text = @'
questions = [
    {
        "title": "You\u2019re a senior database consultant advising a fast-growing tech startup on choosing the right NoSQL database. The company\u2019s application involves complex relationships between data entities and requires efficient traversal between these relationships.Which NoSQL database type should you recommend?",
        "ques_type": 2,
        "options": [
            "Graph database",
            "Key-value store",
            "Document store",
            "Wide-column store"
        ],
        "score": "Graph database"
    },
    {
        "title": "You\u2019re a database administrator at a financial analytics company. You need to improve the performance of a NoSQL database used primarily for read operations. The database is experiencing increased latency due to complex queries on large financial datasets, so you need to improve read performance without impacting data integrity or causing significant downtime.Which action should you take?",
        "ques_type": 2,
        "options": [
            "Optimize the query structure to reduce data scanning.",
            "Increase the overall storage capacity of the database.",
            "Migrate the database to a more powerful server with higher processing capability.",
            "Shift to a relational database model for structured query optimization."
        ],
        "score": "Optimize the query structure to reduce data scanning."
    },
    {
        "title": "You're the lead database architect for a large retail company. Its NoSQL database structure leads to frequent data duplication and inconsistencies, especially in product information. You need to redesign the schema to minimize redundancy while ensuring quick access to the latest product information.Which action should you take?",
        "ques_type": 2,
        "options": [
            "Implement a denormalized schema with reference links to related data.",
            "Normalize all data fully, creating separate collections for each entity.",
            "Keep all data in a single large document for each product.",
            "Use a relational database instead for stricter normalization controls."
        ],
        "score": "Implement a denormalized schema with reference links to related data."
    },
    {
        "title": "You're evaluating NoSQL databases in a distributed systems environment using microservices architecture to ensure efficient data handling and service scalability. The system needs to best support the distributed nature of microservices and ensure high availability and scalability of services.What key feature should you prioritize in the database?",
        "ques_type": 2,
        "options": [
            "Native support for horizontal scaling and data replication.",
            "Single-node database deployment.",
            "ACID properties over BASE principles.",
            "Exclusive use of SQL-like query languages for data access."
        ],
        "score": "Native support for horizontal scaling and data replication."
    }
]
'@
$ws.Range("A2").ClearContents()
$ws.Range("A1").Style = "Normal"
$ws.Range("A1").Value = $text
$ws.Rows(1).AutoFit()
